# Auto-generated edit script applying the Bahamut_Profits scheduled-runner update.
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H:N) for specific
# rows across all 8 job sheets, matching the refreshed market-board snapshot.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 482.68
$ws.Range("I107").Value = 421.05884
$ws.Range("J107").Value = 613.625
$ws.Range("K107").Value = 421.05884
$ws.Range("L107").Value = 613.625
$ws.Range("M107").Value = 1498.94116
$ws.Range("N107").Value = -4453.625
$ws.Range("H137").Value = 901.375
$ws.Range("I137").Value = 816.13635
$ws.Range("J137").Value = 1088.9
$ws.Range("K137").Value = 2448.40905
$ws.Range("L137").Value = 3266.7
$ws.Range("M137").Value = 101.5909499999998
$ws.Range("N137").Value = -8366.700000000001
$ws.Range("H138").Value = 3229.25
$ws.Range("I138").Value = 1108.5428
$ws.Range("K138").Value = 3325.6284
$ws.Range("M138").Value = 1814.3716

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2346.4614
$ws.Range("I2").Value = 1951.5883
$ws.Range("J2").Value = 3092.3333
$ws.Range("K2").Value = 1951.5883
$ws.Range("L2").Value = 3092.3333
$ws.Range("M2").Value = -1838.5883
$ws.Range("N2").Value = -3318.3333
$ws.Range("H32").Value = 14268.118
$ws.Range("I32").Value = 8036.8823
$ws.Range("J32").Value = 80475
$ws.Range("K32").Value = 8036.8823
$ws.Range("L32").Value = 80475
$ws.Range("M32").Value = -7749.8823
$ws.Range("N32").Value = -81049
$ws.Range("H61").Value = 2328.3684
$ws.Range("I61").Value = 2358.4375
$ws.Range("K61").Value = 2358.4375
$ws.Range("M61").Value = -2146.4375
$ws.Range("H97").Value = 2715.625
$ws.Range("I97").Value = 3070.4736
$ws.Range("K97").Value = 3070.4736
$ws.Range("M97").Value = -2574.4736
$ws.Range("H110").Value = 798.1
$ws.Range("I110").Value = 711.5714
$ws.Range("J110").Value = 1000
$ws.Range("K110").Value = 711.5714
$ws.Range("L110").Value = 1000
$ws.Range("M110").Value = 1333.4286
$ws.Range("N110").Value = -5090
$ws.Range("H116").Value = 2346.4614
$ws.Range("I116").Value = 1951.5883
$ws.Range("J116").Value = 3092.3333
$ws.Range("K116").Value = 1951.5883
$ws.Range("L116").Value = 3092.3333
$ws.Range("M116").Value = 342.4117000000001
$ws.Range("N116").Value = -7680.3333
$ws.Range("H136").Value = 2328.3684
$ws.Range("I136").Value = 2358.4375
$ws.Range("K136").Value = 7075.3125
$ws.Range("M136").Value = -4525.3125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2346.4614
$ws.Range("I3").Value = 1951.5883
$ws.Range("J3").Value = 3092.3333
$ws.Range("K3").Value = 1951.5883
$ws.Range("L3").Value = 3092.3333
$ws.Range("M3").Value = -1837.5883
$ws.Range("N3").Value = -3320.3333
$ws.Range("H80").Value = 560.5294
$ws.Range("I80").Value = 401.2
$ws.Range("J80").Value = 626.9167
$ws.Range("K80").Value = 401.2
$ws.Range("L80").Value = 626.9167
$ws.Range("M80").Value = 596.8
$ws.Range("N80").Value = -2622.9167
$ws.Range("H83").Value = 560.5294
$ws.Range("I83").Value = 401.2
$ws.Range("J83").Value = 626.9167
$ws.Range("K83").Value = 2006
$ws.Range("L83").Value = 3134.5835
$ws.Range("M83").Value = 2986
$ws.Range("N83").Value = -13118.5835
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()
$ws.Range("H134").Value = 34522
$ws.Range("I134").Value = 2152.476
$ws.Range("J134").Value = 102498
$ws.Range("K134").Value = 6457.428
$ws.Range("L134").Value = 307494
$ws.Range("M134").Value = -3922.428
$ws.Range("N134").Value = -312564

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H38").Value = 11266.667
$ws.Range("I38").Value = 3800
$ws.Range("J38").Value = 15000
$ws.Range("K38").Value = 3800
$ws.Range("L38").Value = 15000
$ws.Range("M38").Value = -3423
$ws.Range("N38").Value = -15754
$ws.Range("H46").Value = 11266.667
$ws.Range("I46").Value = 3800
$ws.Range("J46").Value = 15000
$ws.Range("K46").Value = 3800
$ws.Range("L46").Value = 15000
$ws.Range("M46").Value = -3589
$ws.Range("N46").Value = -15422
$ws.Range("H134").Value = 23810984
$ws.Range("I134").Value = 1508.7778
$ws.Range("K134").Value = 4526.3334
$ws.Range("M134").Value = -1991.3334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 234
$ws.Range("I13").Value = 251
$ws.Range("K13").Value = 753
$ws.Range("M13").Value = -585
$ws.Range("H106").Value = 8000
$ws.Range("J106").Value = 8000
$ws.Range("L106").Value = 24000
$ws.Range("N106").Value = -25892
$ws.Range("H113").Value = 641.1177
$ws.Range("I113").Value = 669.13336
$ws.Range("J113").Value = 619
$ws.Range("K113").Value = 2007.40008
$ws.Range("L113").Value = 1857
$ws.Range("M113").Value = 162.5999199999999
$ws.Range("N113").Value = -6197
$ws.Range("H131").Value = 18557300
$ws.Range("J131").Value = 1802.2
$ws.Range("L131").Value = 5406.6
$ws.Range("N131").Value = -15486.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1083.6842
$ws.Range("I122").Value = 1005.7143
$ws.Range("J122").Value = 1302
$ws.Range("K122").Value = 3017.1429
$ws.Range("L122").Value = 3906
$ws.Range("M122").Value = -567.1428999999998
$ws.Range("N122").Value = -8806
$ws.Range("H132").Value = 2328.625
$ws.Range("I132").Value = 2081.8572
$ws.Range("J132").Value = 2520.5557
$ws.Range("K132").Value = 6245.571599999999
$ws.Range("L132").Value = 7561.6671
$ws.Range("M132").Value = -3715.571599999999
$ws.Range("N132").Value = -12621.6671

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H18").Value = 22000
$ws.Range("I18").Value = 22000
$ws.Range("K18").Value = 22000
$ws.Range("M18").Value = -21828
$ws.Range("H122").Value = 2174.25
$ws.Range("I122").Value = 2174.25
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 6522.75
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -4072.75
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1359.2307
$ws.Range("I100").Value = 1067
$ws.Range("J100").Value = 2333.3333
$ws.Range("K100").Value = 2134
$ws.Range("L100").Value = 4666.6666
$ws.Range("M100").Value = -1593
$ws.Range("N100").Value = -5748.6666
$ws.Range("H113").Value = 299.85715
$ws.Range("I113").Value = 236
$ws.Range("J113").Value = 385
$ws.Range("K113").Value = 708
$ws.Range("L113").Value = 1155
$ws.Range("M113").Value = 1462
$ws.Range("N113").Value = -5495
$ws.Range("H126").Value = 846.7143
$ws.Range("I126").Value = 646.4
$ws.Range("J126").Value = 1347.5
$ws.Range("K126").Value = 1939.2
$ws.Range("L126").Value = 4042.5
$ws.Range("M126").Value = 530.8000000000002
$ws.Range("N126").Value = -8982.5
$ws.Range("H132").Value = 1702.8182
$ws.Range("I132").Value = 1500.4706
$ws.Range("J132").Value = 2390.8
$ws.Range("K132").Value = 4501.4118
$ws.Range("L132").Value = 7172.400000000001
$ws.Range("M132").Value = -1971.4118
$ws.Range("N132").Value = -12232.4
